# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 = Home stats)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 439
$wsOff.Range("C2").Value = 287
$wsOff.Range("D2").Value = 97
$wsOff.Range("E2").Value = 30

# Update DEF sheet (row 2 = Home stats)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 549
$wsDef.Range("C2").Value = 392
$wsDef.Range("D2").Value = 120
$wsDef.Range("E2").Value = 59
$wsDef.Range("F2").Value = 12
$wsDef.Range("G2").Value = 3
